# Adds a "percentage" column to both sheets of the workbook.
#   - "PI hours" (sheet1): percentage is inserted between "hours" and "dept"
#   - "dept hours" (sheet2): percentage is appended after "hours"
# Percentage = 100 * hours / sum(hours) for each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "PI hours" -- insert a new column D ("percentage") between
# "hours" (C) and "dept" (which shifts from D to E).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a blank column before the current "dept" column (D); this
# shifts "dept" to E and the new column inherits the header formatting.
$ws1.Columns.Item(4).Insert()

# Header cell, reuse the same value/formatting as the other header cells.
$ws1.Cells.Item(1, 4).Value = "percentage"

# Sum up "hours" (column C, rows 2-5) for the percentage base.
$total1 = 0
$lastRow1 = 5
for ($r = 2; $r -le $lastRow1; $r++) {
    $total1 = $total1 + $ws1.Cells.Item($r, 3).Value2
}

for ($r = 2; $r -le $lastRow1; $r++) {
    $hours = $ws1.Cells.Item($r, 3).Value2
    $ws1.Cells.Item($r, 4).Value = 100 * $hours / $total1
}

# Reset page margins to Excel's standard defaults (inches -> points).
$ws1.PageSetup.LeftMargin = 0.7 * 72
$ws1.PageSetup.RightMargin = 0.7 * 72
$ws1.PageSetup.TopMargin = 0.75 * 72
$ws1.PageSetup.BottomMargin = 0.75 * 72
$ws1.PageSetup.HeaderMargin = 0.3 * 72
$ws1.PageSetup.FooterMargin = 0.3 * 72

# ---------------------------------------------------------------------
# Sheet 2: "dept hours" -- append a new column D ("percentage") after
# "hours" (C). "dept" (B) / "hours" (C) keep their positions.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Copy formatting from an existing styled header cell (B1) so the new
# header cell D1 reuses the same cell style instead of creating a new one.
$ws2.Range("B1").Copy($ws2.Range("D1"))
$ws2.Cells.Item(1, 4).Value = "percentage"

$total2 = 0
$lastRow2 = 6
for ($r = 2; $r -le $lastRow2; $r++) {
    $total2 = $total2 + $ws2.Cells.Item($r, 3).Value2
}

for ($r = 2; $r -le $lastRow2; $r++) {
    $hours = $ws2.Cells.Item($r, 3).Value2
    $ws2.Cells.Item($r, 4).Value = 100 * $hours / $total2
}

# Reset page margins to Excel's standard defaults (inches -> points).
$ws2.PageSetup.LeftMargin = 0.7 * 72
$ws2.PageSetup.RightMargin = 0.7 * 72
$ws2.PageSetup.TopMargin = 0.75 * 72
$ws2.PageSetup.BottomMargin = 0.75 * 72
$ws2.PageSetup.HeaderMargin = 0.3 * 72
$ws2.PageSetup.FooterMargin = 0.3 * 72
